$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.374.68"
$ws.Range("E2").Value = "  +0.30%  "

$ws.Range("D3").Value = "'3.074.86"
$ws.Range("E3").Value = "  +0.50%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'514.71"
$ws.Range("E5").Value = "  -0.32%  "

$ws.Range("D6").Value = "'141.34"
$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "'0.434"
$ws.Range("E8").Value = "  -1.57%  "

$ws.Range("E9").Value = "  +0.91%  "

$ws.Range("D10").Value = "'0.108"
$ws.Range("E10").Value = "  -1.16%  "

$ws.Range("D11").Value = "'0.372"
$ws.Range("E11").Value = "  -1.43%  "

$ws.Range("D12").Value = "'3.612.10"
$ws.Range("E12").Value = "  +0.97%  "

$ws.Range("E13").Value = "  +2.78%  "

$ws.Range("D14").Value = "'25.51"
$ws.Range("E14").Value = "  -5.08%  "

$ws.Range("E15").Value = "  -1.59%  "

$ws.Range("D16").Value = "'57.505.35"
$ws.Range("E16").Value = "  +0.66%  "

$ws.Range("D17").Value = "'3.085.13"
$ws.Range("E17").Value = "  +0.98%  "

$ws.Range("D18").Value = "'6.08"
$ws.Range("E18").Value = "  -1.12%  "

$ws.Range("D19").Value = "'13.16"
$ws.Range("E19").Value = "  -1.99%  "

$ws.Range("D20").Value = "'8.12"
$ws.Range("E20").Value = "  +0.02%  "

$ws.Range("D21").Value = "'332.72"
$ws.Range("E21").Value = "  +0.33%  "

$ws.Range("E22").Value = "  -0.31%  "

$ws.Range("D23").Value = "'0.500"
$ws.Range("E23").Value = "  -1.52%  "

$ws.Range("D24").Value = "'65.71"
$ws.Range("E24").Value = "  -0.21%  "

$ws.Range("D25").Value = "'0.170"
$ws.Range("E25").Value = "  +3.30%  "

$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("D27").Value = "'0.0₃0910"
$ws.Range("E27").Value = "  +1.87%  "

$ws.Range("D28").Value = "'6.36"
$ws.Range("E28").Value = "  -5.61%  "

$ws.Range("D29").Value = "'7.17"
$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("D30").Value = "'1.81"
$ws.Range("E30").Value = "  -0.14%  "

$ws.Range("D31").Value = "'20.79"
$ws.Range("E31").Value = "  -0.06%  "

$ws.Range("E32").Value = "  -3.64%  "

$ws.Range("D33").Value = "'154.77"
$ws.Range("E33").Value = "  +1.89%  "

$ws.Range("D34").Value = "'27.79"
$ws.Range("E34").Value = "  +10.29%  "

$ws.Range("D35").Value = "'4.53"
$ws.Range("E35").Value = "  -3.69%  "

$ws.Range("D36").Value = "'5.84"
$ws.Range("E36").Value = "  -1.72%  "

$ws.Range("D37").Value = "'1.25"
$ws.Range("E37").Value = "  -1.33%  "

$ws.Range("D38").Value = "'0.0671"
$ws.Range("E38").Value = "  -0.21%  "

$ws.Range("D39").Value = "'3.121.54"
$ws.Range("E39").Value = "  +0.93%  "

$ws.Range("D40").Value = "'36.82"
$ws.Range("E40").Value = "  -0.86%  "

$ws.Range("D41").Value = "'3.87"
$ws.Range("E41").Value = "  -0.27%  "

$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.08%  "

$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "'0.665"
$ws.Range("E43").Value = "  +0.36%  "

$ws.Range("D44").Value = "'2.280.61"
$ws.Range("E44").Value = "  +3.66%  "

$ws.Range("D45").Value = "'0.0258"
$ws.Range("E45").Value = "  +6.90%  "

$ws.Range("E46").Value = "  -0.98%  "

$ws.Range("E47").Value = "  -2.26%  "

$ws.Range("D48").Value = "'19.97"
$ws.Range("E48").Value = "  -2.06%  "

$ws.Range("D49").Value = "'5.87"
$ws.Range("E49").Value = "  -2.96%  "

$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "'0.0874"
$ws.Range("E50").Value = "  +0.71%  "

$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D51").Value = "'249.90"
$ws.Range("E51").Value = "  +5.10%  "
